# Refresh the cryptos list snapshot (price + 1h volume change) and fix
# the OKB / Stacks row ordering that had swapped.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range('D2').Value = '60.991.61'
$ws.Range('E2').Value = '  -3.02%  '

# Row 3: Ethereum
$ws.Range('D3').Value = '2.996.88'

# Row 4: TetherUSD
$ws.Range('D4').Value = '''1.00'
$ws.Range('E4').Value = '  +0.19%  '

# Row 5: BNB
$ws.Range('D5').Value = '''533.97'
$ws.Range('E5').Value = '  -0.27%  '

# Row 6: Solana
$ws.Range('E6').Value = '  +1.10%  '

# Row 7: USDC
$ws.Range('E7').Value = '  +0.02%  '

# Row 8: LidoStakedEther
$ws.Range('D8').Value = '2.996.99'
$ws.Range('E8').Value = '  -1.83%  '

# Row 9: XRP
$ws.Range('D9').Value = '''0.497'
$ws.Range('E9').Value = '  +1.32%  '

# Row 10: Dogecoin
$ws.Range('E10').Value = '  -3.23%  '

# Row 11: Toncoin
$ws.Range('D11').Value = '''6.08'
$ws.Range('E11').Value = '  -0.23%  '

# Row 12: Cardano
$ws.Range('E12').Value = '  -0.71%  '

# Row 13: ShibaInu
$ws.Range('E13').Value = '  -0.99%  '

# Row 14: Avalanche
$ws.Range('D14').Value = '''34.21'
$ws.Range('E14').Value = '  +0.65%  '

# Row 15: WrappedliquidstakedEther2.0
$ws.Range('D15').Value = '3.491.64'
$ws.Range('E15').Value = '  -1.83%  '

# Row 16: TRON
$ws.Range('E16').Value = '  -0.30%  '

# Row 17: WrappedBTC
$ws.Range('D17').Value = '61.179.29'
$ws.Range('E17').Value = '  -2.69%  '

# Row 18: WrappedEther
$ws.Range('D18').Value = '3.002.33'
$ws.Range('E18').Value = '  -2.10%  '

# Row 19: Polkadot
$ws.Range('E19').Value = '  +0.47%  '

# Row 20: BitcoinCash
$ws.Range('D20').Value = '''463.94'
$ws.Range('E20').Value = '  -3.15%  '

# Row 21: Chainlink
$ws.Range('D21').Value = '''13.21'
$ws.Range('E21').Value = '  -0.12%  '

# Row 22: Polygon
$ws.Range('E22').Value = '  -2.15%  '

# Row 23: Uniswap
$ws.Range('E23').Value = '  -1.46%  '

# Row 24: Litecoin
$ws.Range('D24').Value = '''79.17'
$ws.Range('E24').Value = '  +0.77%  '

# Row 25: InternetComputer(DFINITY)
$ws.Range('D25').Value = '''12.06'
$ws.Range('E25').Value = '  +0.51%  '

# Row 26: Dai
$ws.Range('D26').Value = '''1.00'
$ws.Range('E26').Value = '  +0.22%  '

# Row 27: PancakeSwap
$ws.Range('E27').Value = '  -0.56%  '

# Row 28: RenderToken
$ws.Range('E28').Value = '  -2.12%  '

# Row 29: FirstDigitalUSD
$ws.Range('D29').Value = '''0.998'
$ws.Range('E29').Value = '  -0.20%  '

# Row 30: ImmutableX
$ws.Range('D30').Value = '''1.90'
$ws.Range('E30').Value = '  +1.76%  '

# Row 31: EthereumClassic
$ws.Range('D31').Value = '''25.46'
$ws.Range('E31').Value = '  -1.33%  '

# Row 32: Mantle
$ws.Range('E32').Value = '  +3.66%  '

# Row 33: NEARProtocol
$ws.Range('E33').Value = '  +3.04%  '

# Row 34: Stacks
$ws.Range('B34').Value = 'Stacks'
$ws.Range('C34').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D34').Value = '''2.28'
$ws.Range('E34').Value = '  -2.57%  '

# Row 35: OKB
$ws.Range('B35').Value = 'OKB'
$ws.Range('C35').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D35').Value = '''55.38'
$ws.Range('E35').Value = '  -2.65%  '

# Row 36: Filecoin
$ws.Range('D36').Value = '''5.88'
$ws.Range('E36').Value = '  -1.26%  '

# Row 37: Bittensor
$ws.Range('D37').Value = '''457.68'
$ws.Range('E37').Value = '  -3.54%  '

# Row 38: Maker
$ws.Range('D38').Value = '3.213.32'
$ws.Range('E38').Value = '  +4.13%  '

# Row 39: Hedera
$ws.Range('D39').Value = '''0.0787'
$ws.Range('E39').Value = '  -0.11%  '

# Row 40: VeChain
$ws.Range('D40').Value = '''0.0384'
$ws.Range('E40').Value = '  -1.73%  '

# Row 41: Kaspa
$ws.Range('E41').Value = '  +2.85%  '

# Row 42: Cosmos
$ws.Range('E42').Value = '  +1.53%  '

# Row 43: InjectiveProtocol
$ws.Range('D43').Value = '''27.81'
$ws.Range('E43').Value = '  +15.23%  '

# Row 44: dogwifhat
$ws.Range('E44').Value = '  -5.23%  '

# Row 46: TheGraph
$ws.Range('E46').Value = '  -1.50%  '

# Row 47: Fetch.AI
$ws.Range('E47').Value = '  +0.55%  '

# Row 48: Monero
$ws.Range('D48').Value = '''119.70'
$ws.Range('E48').Value = '  -1.31%  '

# Row 49: Stellar
$ws.Range('E49').Value = '  +1.00%  '

# Row 50: PEPE
$ws.Range('D50').Value = '0.0₃0493'
$ws.Range('E50').Value = '  -7.22%  '

# Row 51: BitgetToken
$ws.Range('E51').Value = '  +8.51%  '
